# Bai 29 - DataProvider TestNG
$wb = $excel.ActiveWorkbook

# Rename Sheet1 -> LoginSuccess
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "LoginSuccess"

# Populate the LoginSuccess sheet with header + data rows (emails copied from Login sheet, password 123456)
$ws2.Range("A1").Value = "EMAIL"
$ws2.Range("B1").Value = "PASSWORD"

$emails = @("admin@example.com", "admin123@example.com", "admin@example.com", "admin@example.com", "admin@example.com", "admin123@example.com", "admin@example.com")

for ($i = 0; $i -lt $emails.Length; $i++) {
    $row = $i + 2
    $ws2.Range("A$row").Value = $emails[$i]
    $ws2.Range("B$row").Value = "'123456"
    $ws2.Hyperlinks.Add($ws2.Range("A$row"), "mailto:" + $emails[$i])
}

$ws2.Columns.Item(1).ColumnWidth = 24.3828125
$ws2.Columns.Item(2).ColumnWidth = 15.69140625

# Login sheet keeps its own remembered selection (B11), but LoginSuccess
# becomes the active sheet/tab with its own selection (B10) - select it last.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B11").Select() | Out-Null

$ws2.Select() | Out-Null
$ws2.Range("B10").Select() | Out-Null
